$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text-safe cell updates (values that will not be mis-typed as numbers) ---
$ws.Range("D2").Value = '96.971.23'
$ws.Range("E2").Value = '  -0.16%  '
$ws.Range("D3").Value = '3.675.73'
$ws.Range("E3").Value = '  +2.86%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("E6").Value = '  +10.19%  '
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("E9").Value = '  +3.63%  '
$ws.Range("D11").Value = '3.675.03'
$ws.Range("E11").Value = '  +2.89%  '
$ws.Range("E12").Value = '  +3.29%  '
$ws.Range("E13").Value = '  +1.27%  '
$ws.Range("E14").Value = '  +7.27%  '
$ws.Range("D15").Value = '4.360.65'
$ws.Range("E15").Value = '  +2.92%  '
$ws.Range("E16").Value = '  +3.58%  '
$ws.Range("D17").Value = '96.802.86'
$ws.Range("E17").Value = '  -0.12%  '
$ws.Range("E18").Value = '  +4.75%  '
$ws.Range("D19").Value = '3.681.82'
$ws.Range("E19").Value = '  +3.17%  '
$ws.Range("E20").Value = '  +6.08%  '
$ws.Range("E21").Value = '  +0.54%  '
$ws.Range("E22").Value = '  +0.81%  '
$ws.Range("E23").Value = '  +3.57%  '
$ws.Range("E24").Value = '  +0.58%  '
$ws.Range("E25").Value = '  +4.45%  '
$ws.Range("E26").Value = '  -0.57%  '
$ws.Range("E27").Value = '  +0.98%  '
$ws.Range("E28").Value = '  +3.48%  '
$ws.Range("D29").Value = '3.870.25'
$ws.Range("E29").Value = '  +2.81%  '
$ws.Range("E30").Value = '  +0.93%  '
$ws.Range("E31").Value = '  +5.43%  '
$ws.Range("E32").Value = '  +1.93%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("E34").Value = '  +15.20%  '
$ws.Range("E35").Value = '  +1.33%  '
$ws.Range("E38").Value = '  +6.46%  '
$ws.Range("E39").Value = '  +7.83%  '
$ws.Range("E40").Value = '  +1.59%  '
$ws.Range("E41").Value = '  +16.17%  '
$ws.Range("E42").Value = '  +5.48%  '
$ws.Range("E43").Value = '  +2.43%  '
$ws.Range("E44").Value = '  +4.67%  '
$ws.Range("E45").Value = '  +16.11%  '
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("E48").Value = '  +5.21%  '
$ws.Range("E49").Value = '  +1.17%  '
$ws.Range("E50").Value = '  +3.29%  '
$ws.Range("E51").Value = '  +0.18%  '

# --- Numeric-looking values that must be forced to remain Text ---
$numericTextCells = @("D5", "D7", "D8", "D9", "D12", "D18", "D20", "D22", "D23", "D25", "D27", "D28", "D31", "D33", "D34", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D50")
foreach ($cellRef in $numericTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}
$ws.Range("D5").Value = '239.46'
$ws.Range("D7").Value = '653.59'
$ws.Range("D8").Value = '0.424'
$ws.Range("D9").Value = '1.09'
$ws.Range("D12").Value = '45.58'
$ws.Range("D18").Value = '9.04'
$ws.Range("D20").Value = '19.07'
$ws.Range("D22").Value = '0.532'
$ws.Range("D23").Value = '531.63'
$ws.Range("D25").Value = '7.16'
$ws.Range("D27").Value = '102.62'
$ws.Range("D28").Value = '13.51'
$ws.Range("D31").Value = '12.52'
$ws.Range("D33").Value = '0.999'
$ws.Range("D34").Value = '1.90'
$ws.Range("D38").Value = '656.51'
$ws.Range("D39").Value = '0.608'
$ws.Range("D40").Value = '8.99'
$ws.Range("D41").Value = '6.99'
$ws.Range("D42").Value = '0.163'
$ws.Range("D43").Value = '2.01'
$ws.Range("D44").Value = '0.965'
$ws.Range("D45").Value = '38.30'
$ws.Range("D47").Value = '0.451'
$ws.Range("D48").Value = '0.0461'
$ws.Range("D50").Value = '8.78'

# --- Row 36 / Row 37 swap (EthereumClassic <-> Binance-PegBSC-USD) plus updated values ---
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.01"
$ws.Range("E36").Value = '  +0.65%  '

$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "32.69"
$ws.Range("E37").Value = '  +2.93%  '
